# Updated symbol list on Sat Dec 31 01:53:27 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# --- Rows 2-9: price (D) updates, Hora (G) flips to 1 ---
Set-TextCell "D2" "245.68"
Set-TextCell "D3" "25.47"
Set-TextCell "D4" "5.084"
Set-TextCell "D5" "0.05600"
Set-TextCell "D6" "6.547"
Set-TextCell "D7" "3.012"
Set-TextCell "D8" "0.8190"
Set-TextCell "D9" "0.8384"

# --- Rows 10-16: coin list rotated up by one (coin that was in row N+1
#     moves to row N), with freshly-updated prices ---
Set-TextCell "B10" "WazirX"
Set-TextCell "C10" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextCell "D10" "0.1341"
Set-TextCell "E10" "9WazirXWRX"

Set-TextCell "B11" "MandalaExchangeToken"
Set-TextCell "C11" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextCell "D11" "0.06927"
Set-TextCell "E11" "10MandalaExchangeTokenMDX"

Set-TextCell "B12" "LiechtensteinCryptoassetsExchange"
Set-TextCell "C12" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextCell "D12" "0.03195"
Set-TextCell "E12" "11LiechtensteinCryptoassetsExchangeLCX"

Set-TextCell "B13" "BitrueCoin"
Set-TextCell "C13" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextCell "D13" "0.02834"
Set-TextCell "E13" "12BitrueCoinBTR"

Set-TextCell "B14" "BitMartToken"
Set-TextCell "C14" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextCell "D14" "0.09392"
Set-TextCell "E14" "13BitMartTokenBMX"

Set-TextCell "B15" "BitForexToken"
Set-TextCell "C15" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextCell "D15" "0.001520"
Set-TextCell "E15" "14BitForexTokenBF"

Set-TextCell "B16" "One"
Set-TextCell "C16" "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextCell "D16" "0.0005962"
Set-TextCell "E16" "15OneONEWorstin24h"

# --- Rows 17-51: remaining price (D) and label (E) updates ---
Set-TextCell "D17" "0.006250"
Set-TextCell "D18" "3.524"
Set-TextCell "D19" "2.073"
Set-TextCell "D21" "0.1338"
Set-TextCell "D22" "3.743"
Set-TextCell "D23" "0.04707"
Set-TextCell "D25" "0.001243"
Set-TextCell "D26" "0.004276"
Set-TextCell "D27" "0.00009705"
Set-TextCell "E27" "26NitroExNTXBestin24h"
Set-TextCell "D28" "0.0001941"
Set-TextCell "D40" "0.03663"
Set-TextCell "D41" "0.006293"
Set-TextCell "E41" "40KickTokenKICK"
Set-TextCell "D42" "0.1054"
Set-TextCell "D43" "0.002679"
Set-TextCell "D44" "0.008386"
Set-TextCell "D45" "0.00005298"
Set-TextCell "D47" "0.2251"
Set-TextCell "E47" "46CoinbaseStockTokenCOIN"
Set-TextCell "D48" "0.002285"
Set-TextCell "D49" "0.00002101"
Set-TextCell "D50" "0.0002001"

# --- Rows 2-51: Hora (G) flips from 0 to 1 across the board ---
for ($row = 2; $row -le 51; $row++) {
    Set-TextCell "G$row" "1"
}
